$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 999.3333
$ws.Range("I12").Value = 99.5
$ws.Range("J12").Value = 2799
$ws.Range("K12").Value = 99.5
$ws.Range("L12").Value = 2799
$ws.Range("M12").Value = 70.5
$ws.Range("N12").Value = -3139

$ws.Range("H51").Value = 2870.8
$ws.Range("J51").Value = 3550.8333
$ws.Range("L51").Value = 3550.8333
$ws.Range("N51").Value = -4518.8333

$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 150000
$ws.Range("N57").Value = -150998

$ws.Range("H97").Value = 47726.57
$ws.Range("I97").Value = 5450
$ws.Range("J97").Value = 64637.2
$ws.Range("K97").Value = 16350
$ws.Range("L97").Value = 193911.6
$ws.Range("M97").Value = -15854
$ws.Range("N97").Value = -194903.6

$ws.Range("H100").Value = 1670.6111
$ws.Range("I100").Value = 1180.6471
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 1180.6471
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -639.6470999999999
$ws.Range("N100").Value = -11082

$ws.Range("H103").Value = 463.2857
$ws.Range("J103").Value = 433
$ws.Range("L103").Value = 1299
$ws.Range("N103").Value = -2471

$ws.Range("H134").Value = 115833
$ws.Range("J134").Value = 115833
$ws.Range("L134").Value = 115833
$ws.Range("N134").Value = -125973

$ws.Range("H138").Value = 5323654.5
$ws.Range("J138").Value = 6254892.5
$ws.Range("L138").Value = 18764677.5
$ws.Range("N138").Value = -18774957.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 969.8570999999999
$ws.Range("I4").Value = 247.5
$ws.Range("J4").Value = 1933
$ws.Range("K4").Value = 247.5
$ws.Range("L4").Value = 1933
$ws.Range("M4").Value = -131.5
$ws.Range("N4").Value = -2165

$ws.Range("H32").Value = 12018.116
$ws.Range("I32").Value = 6819.1353
$ws.Range("K32").Value = 6819.1353
$ws.Range("M32").Value = -6532.1353

$ws.Range("H61").Value = 184552.61
$ws.Range("I61").Value = 1961.2
$ws.Range("K61").Value = 1961.2
$ws.Range("M61").Value = -1749.2

$ws.Range("H102").Value = 4959.0835
$ws.Range("I102").Value = 4955.364
$ws.Range("K102").Value = 4955.364
$ws.Range("M102").Value = -3333.364

$ws.Range("H110").Value = 7738.7144
$ws.Range("I110").Value = 11161.5
$ws.Range("J110").Value = 4627.091
$ws.Range("K110").Value = 11161.5
$ws.Range("L110").Value = 4627.091
$ws.Range("M110").Value = -9116.5
$ws.Range("N110").Value = -8717.091

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 3206.25
$ws.Range("I122").Value = 2504
$ws.Range("K122").Value = 7512
$ws.Range("M122").Value = -5062

$ws.Range("H136").Value = 184552.61
$ws.Range("I136").Value = 1961.2
$ws.Range("K136").Value = 5883.6
$ws.Range("M136").Value = -3333.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6260
$ws.Range("I86").Value = 6374.8184
$ws.Range("J86").Value = 4997
$ws.Range("K86").Value = 6374.8184
$ws.Range("L86").Value = 4997
$ws.Range("M86").Value = -5251.8184
$ws.Range("N86").Value = -7243

$ws.Range("H89").Value = 6260
$ws.Range("I89").Value = 6374.8184
$ws.Range("J89").Value = 4997
$ws.Range("K89").Value = 31874.092
$ws.Range("L89").Value = 24985
$ws.Range("M89").Value = -26258.092
$ws.Range("N89").Value = -36217

$ws.Range("H94").Value = 1889.7858
$ws.Range("I94").Value = 1889.7858
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1889.7858
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1438.7858
$ws.Range("N94").ClearContents()

$ws.Range("H97").Value = 8508.454
$ws.Range("J97").Value = 29998.5
$ws.Range("L97").Value = 29998.5
$ws.Range("N97").Value = -31980.5

$ws.Range("H99").Value = 5690.1113
$ws.Range("I99").Value = 3552.3
$ws.Range("K99").Value = 3552.3
$ws.Range("M99").Value = -2054.3

$ws.Range("H105").Value = 7275.3335
$ws.Range("I105").Value = 7530.4
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 7530.4
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -5783.4
$ws.Range("N105").Value = -9494

$ws.Range("H134").Value = 1245
$ws.Range("I134").Value = 1188.3871
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 3565.1613
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1030.1613
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 708.6875
$ws.Range("I7").Value = 1061.9
$ws.Range("K7").Value = 1061.9
$ws.Range("M7").Value = -948.9000000000001

$ws.Range("H22").Value = 509.625
$ws.Range("I22").Value = 414.6
$ws.Range("K22").Value = 414.6
$ws.Range("M22").Value = -64.60000000000002

$ws.Range("H31").Value = 2902.9443
$ws.Range("I31").Value = 2053.2
$ws.Range("K31").Value = 2053.2
$ws.Range("M31").Value = -1758.2

$ws.Range("H34").Value = 2902.9443
$ws.Range("I34").Value = 2053.2
$ws.Range("K34").Value = 2053.2
$ws.Range("M34").Value = -1851.2

$ws.Range("H99").Value = 3144.3914
$ws.Range("I99").Value = 2706.5625
$ws.Range("K99").Value = 2706.5625
$ws.Range("M99").Value = -1208.5625

$ws.Range("H105").Value = 2130.2307
$ws.Range("I105").Value = 2058.2
$ws.Range("K105").Value = 2058.2
$ws.Range("M105").Value = -311.1999999999998

$ws.Range("H126").Value = 3144.3914
$ws.Range("I126").Value = 2706.5625
$ws.Range("K126").Value = 8119.6875
$ws.Range("M126").Value = -5649.6875

$ws.Range("H132").Value = 1402.9773
$ws.Range("I132").Value = 1265.0256
$ws.Range("K132").Value = 3795.0768
$ws.Range("M132").Value = -1265.0768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 19.5
$ws.Range("I2").Value = 19.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 117
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4
$ws.Range("N2").ClearContents()

$ws.Range("H12").Value = 401.8125
$ws.Range("I12").Value = 183.6
$ws.Range("J12").Value = 501
$ws.Range("K12").Value = 550.8
$ws.Range("L12").Value = 1503
$ws.Range("M12").Value = -377.8
$ws.Range("N12").Value = -1849

$ws.Range("H137").Value = 73704.64
$ws.Range("I137").Value = 92896.82000000001
$ws.Range("J137").Value = 3333.3333
$ws.Range("K137").Value = 278690.46
$ws.Range("L137").Value = 9999.999899999999
$ws.Range("M137").Value = -273590.46
$ws.Range("N137").Value = -20199.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6886.0835
$ws.Range("I80").Value = 3128
$ws.Range("J80").Value = 9570.429
$ws.Range("K80").Value = 3128
$ws.Range("L80").Value = 9570.429
$ws.Range("M80").Value = -2130
$ws.Range("N80").Value = -11566.429

$ws.Range("H83").Value = 6886.0835
$ws.Range("I83").Value = 3128
$ws.Range("J83").Value = 9570.429
$ws.Range("K83").Value = 15640
$ws.Range("L83").Value = 47852.145
$ws.Range("M83").Value = -10648
$ws.Range("N83").Value = -57836.145

$ws.Range("H132").Value = 3125.4092
$ws.Range("J132").Value = 4987.375
$ws.Range("L132").Value = 14962.125
$ws.Range("N132").Value = -20022.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4224.1025
$ws.Range("I40").Value = 3530.4517
$ws.Range("J40").Value = 6912
$ws.Range("K40").Value = 3530.4517
$ws.Range("L40").Value = 6912
$ws.Range("M40").Value = -3394.4517
$ws.Range("N40").Value = -7184

$ws.Range("H46").Value = 3307.1667
$ws.Range("I46").Value = 2233.6667
$ws.Range("J46").Value = 3521.8667
$ws.Range("K46").Value = 2233.6667
$ws.Range("L46").Value = 3521.8667
$ws.Range("M46").Value = -2045.6667
$ws.Range("N46").Value = -3897.8667

$ws.Range("H100").Value = 1002200
$ws.Range("I100").Value = 4400
$ws.Range("K100").Value = 4400
$ws.Range("M100").Value = -3859

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws.Range("H132").Value = 2839.5757
$ws.Range("I132").Value = 2160.625
$ws.Range("J132").Value = 4650.1113
$ws.Range("K132").Value = 6481.875
$ws.Range("L132").Value = 13950.3339
$ws.Range("M132").Value = -3951.875
$ws.Range("N132").Value = -19010.3339

$ws.Range("H136").Value = 4450.8276
$ws.Range("I136").Value = 3571.3635
$ws.Range("J136").Value = 7214.857
$ws.Range("K136").Value = 10714.0905
$ws.Range("L136").Value = 21644.571
$ws.Range("M136").Value = -8164.0905
$ws.Range("N136").Value = -26744.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 12995.2
$ws.Range("J63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16248

$ws.Range("H66").Value = 12995.2
$ws.Range("J66").Value = 15000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -51240

$ws.Range("H113").Value = 358359.28
$ws.Range("I113").Value = 1322.1428
$ws.Range("K113").Value = 3966.4284
$ws.Range("M113").Value = -1796.4284

$ws.Range("H132").Value = 1919.9016
$ws.Range("I132").Value = 1781.2181
$ws.Range("K132").Value = 5343.6543
$ws.Range("M132").Value = -2813.6543
